# Update workbook with new daily data rows (through 26/04/2021, serial 44312)
# Mirrors the pattern established by existing rows: column A holds the date
# (serial number, formatted via the existing date style), columns B:AX hold
# per-comune counts for that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 233
$firstNewRow = 234
$lastNewRow = 238

# Extend formatting (including the bordered/centered date style in column A)
# from the last existing row down across the new rows, matching how the
# sheet already looks for every prior day.
$ws.Range("A$lastRow`:AX$lastRow").Copy()
$ws.Range("A$firstNewRow`:AX$lastNewRow").PasteSpecial(-4122)

$newData = @(
    @(44308,4,1,2,5,4,2,2,1,0,0,1,5,9,1,0,0,8,1,1,4,55,2,1,1,8,1,0,1,0,10,0,18,2,3,1,2,3,5,3,7,179,1,0,0,0,3,0,1,0),
    @(44309,2,1,0,11,18,0,2,2,1,1,2,3,6,0,0,0,14,1,0,5,32,1,4,8,8,0,0,1,2,3,1,22,1,2,0,5,1,12,4,7,187,2,0,0,0,0,1,1,0),
    @(44310,4,2,1,14,6,0,7,0,0,0,0,9,7,0,0,0,6,1,2,5,39,1,1,2,9,0,5,0,2,3,0,15,1,1,0,2,1,10,0,2,160,1,0,0,0,0,0,1,0),
    @(44311,5,2,5,14,9,3,1,1,3,0,3,2,3,0,1,0,1,1,3,4,39,2,3,1,10,0,0,0,0,0,2,2,3,0,0,3,5,17,0,3,157,6,0,0,0,0,0,0,0),
    @(44312,7,0,0,7,7,2,2,1,6,0,1,5,8,0,0,0,7,0,1,4,61,3,2,0,3,0,0,0,1,2,1,13,2,0,0,0,5,5,0,4,161,0,0,1,0,0,0,0,0)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowValues = $newData[$i]
    $targetRow = $firstNewRow + $i
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($targetRow, $col + 1).Value2 = $rowValues[$col]
    }
}
